# Apply the "Deploying to gh-pages" content update to the FHIR
# StructureDefinition metadata/elements workbook.
#
# Sheet 1 ("Metadata") holds a Property/Value table:
#   - Version bump 5.0.0 -> 6.0.0
#   - Date bump to the new publish timestamp
#   - Publisher filled in ("Alvearie Team")
#   - The duplicated "Contact" / "No display for ContactDetail" rows are
#     replaced by a single "Jurisdiction" / "United States of America" row
#
# Sheet 2 ("Elements") holds the per-element definition table; the root
# Extension element's Short/Definition text is refreshed to match the
# profile's own Title/Description.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsElements = $wb.Worksheets.Item(2)

# --- Sheet 1: Metadata ------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date: refresh publish timestamp
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher: was blank
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -- repurpose it
# as the new "Jurisdiction" / "United States of America" row.
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# remove it entirely so everything below shifts up by one row.
$wsMeta.Rows.Item(11).Delete()

# --- Sheet 2: Elements -------------------------------------------------

# Root Extension element (row 2): Short + Definition now mirror the
# StructureDefinition's own Title and Description instead of the generic
# "Extension" / "An Extension" placeholders.
$wsElements.Range("K2").Value = "Enrollment PCP Name On Enrollment"
$wsElements.Range("L2").Value = "Original primary care physician (PCP) name as reporting on the eligibility record"
